$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$newRow = 68

# Columns A-D hold values that look like dates/times/numbers (e.g.
# "2025-02-13", "23:22:14", "06") but must stay literal text, matching every
# other row already in the sheet. Force a temporary text number format so
# Excel does not auto-convert them to date/time/numeric values, assign the
# literal strings, then clear the temporary formatting again so the new row
# ends up with no explicit per-cell style - exactly like the existing rows.
$textRange = $ws.Range("A" + $newRow + ":D" + $newRow)
$textRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-02-13"
$ws.Cells.Item($newRow, 2).Value = "23:22:14"
$ws.Cells.Item($newRow, 3).Value = "Thursday"
$ws.Cells.Item($newRow, 4).Value = "06"

$textRange.ClearFormats()

# Columns E-T are plain numeric resale figures.
$ws.Cells.Item($newRow, 5).Value = 120745
$ws.Cells.Item($newRow, 6).Value = 142319
$ws.Cells.Item($newRow, 7).Value = 170123
$ws.Cells.Item($newRow, 8).Value = 159044
$ws.Cells.Item($newRow, 9).Value = -1
$ws.Cells.Item($newRow, 10).Value = 144876
$ws.Cells.Item($newRow, 11).Value = -1
$ws.Cells.Item($newRow, 12).Value = -1
$ws.Cells.Item($newRow, 13).Value = 192065
$ws.Cells.Item($newRow, 14).Value = 115196
$ws.Cells.Item($newRow, 15).Value = 45048
$ws.Cells.Item($newRow, 16).Value = 28636
$ws.Cells.Item($newRow, 17).Value = 65864
$ws.Cells.Item($newRow, 18).Value = -1
$ws.Cells.Item($newRow, 19).Value = 45370
$ws.Cells.Item($newRow, 20).Value = -1

$wb.Save()
